$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.530.93"
$ws.Range("E2").Value = "'  +0.40%  "
$ws.Range("D3").Value = "'2.488.21"
$ws.Range("E3").Value = "'  +0.89%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "'  -0.30%  "
$ws.Range("D5").Value = "'313.06"
$ws.Range("E5").Value = "'  +0.19%  "
$ws.Range("D6").Value = "'93.23"
$ws.Range("E6").Value = "'  -1.34%  "
$ws.Range("D7").Value = "'0.546"
$ws.Range("E7").Value = "'  -1.11%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "'  -0.25%  "
$ws.Range("E9").Value = "'  -0.85%  "
$ws.Range("E10").Value = "'  -3.10%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "'  +0.32%  "
$ws.Range("E12").Value = "'  +1.77%  "
$ws.Range("D13").Value = "'2.872.86"
$ws.Range("E13").Value = "'  +0.81%  "
$ws.Range("D14").Value = "'6.85"
$ws.Range("E14").Value = "'  -2.70%  "
$ws.Range("D15").Value = "'15.85"
$ws.Range("E15").Value = "'  +8.48%  "
$ws.Range("D16").Value = "'2.523.60"
$ws.Range("E16").Value = "'  +2.85%  "
$ws.Range("D17").Value = "'0.754"
$ws.Range("E17").Value = "'  -4.24%  "
$ws.Range("D18").Value = "'41.575.54"
$ws.Range("E18").Value = "'  +0.61%  "
$ws.Range("E19").Value = "'  +0.14%  "
$ws.Range("D20").Value = "'0.0₃0928"
$ws.Range("E20").Value = "'  +0.99%  "
$ws.Range("D21").Value = "'71.46"
$ws.Range("E21").Value = "'  +4.97%  "
$ws.Range("D22").Value = "'11.24"
$ws.Range("E22").Value = "'  -2.23%  "
$ws.Range("D23").Value = "'235.80"
$ws.Range("E23").Value = "'  -0.61%  "
$ws.Range("D24").Value = "'2.71"
$ws.Range("E24").Value = "'  -2.73%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "'  +0.01%  "
$ws.Range("D26").Value = "'1.91"
$ws.Range("E26").Value = "'  -1.33%  "
$ws.Range("D27").Value = "'24.99"
$ws.Range("E27").Value = "'  +2.19%  "
$ws.Range("E28").Value = "'  -0.24%  "
$ws.Range("D29").Value = "'9.67"
$ws.Range("E29").Value = "'  -0.48%  "
$ws.Range("D30").Value = "'36.17"
$ws.Range("E30").Value = "'  +0.05%  "
$ws.Range("D31").Value = "'157.39"
$ws.Range("E31").Value = "'  +2.65%  "
$ws.Range("D32").Value = "'5.45"
$ws.Range("E32").Value = "'  -2.64%  "
$ws.Range("E33").Value = "'  -1.36%  "
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0756"
$ws.Range("E34").Value = "'  +0.04%  "
$ws.Range("B35").Value = "'Celestia"
$ws.Range("C35").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'17.99"
$ws.Range("E35").Value = "'  +6.06%  "
$ws.Range("D36").Value = "'2.45"
$ws.Range("E36").Value = "'  -6.05%  "
$ws.Range("D37").Value = "'2.94"
$ws.Range("E37").Value = "'  -2.43%  "
$ws.Range("E38").Value = "'  +2.23%  "
$ws.Range("D39").Value = "'1.85"
$ws.Range("E39").Value = "'  -2.43%  "
$ws.Range("E40").Value = "'  -0.13%  "
$ws.Range("D41").Value = "'4.16"
$ws.Range("E41").Value = "'  -2.70%  "
$ws.Range("E42").Value = "'  -0.34%  "
$ws.Range("D43").Value = "'19.98"
$ws.Range("E43").Value = "'  -5.85%  "
$ws.Range("D44").Value = "'1.970.18"
$ws.Range("E44").Value = "'  -0.48%  "
$ws.Range("E45").Value = "'  -0.65%  "
$ws.Range("E46").Value = "'  -3.31%  "
$ws.Range("E47").Value = "'  +2.01%  "
$ws.Range("D48").Value = "'2.728.60"
$ws.Range("E48").Value = "'  +0.67%  "
$ws.Range("D49").Value = "'96.64"
$ws.Range("E49").Value = "'  -0.45%  "
$ws.Range("D50").Value = "'67.85"
$ws.Range("E50").Value = "'  -3.03%  "
$ws.Range("D51").Value = "'73.72"
$ws.Range("E51").Value = "'  -3.09%  "
